# Auto-generated Excel COM-interop script
# Applies the scheduled-runner market-price refresh described in the commit diff
# to each affected sheet/row/cell. No formulas are involved -> plain numeric writes.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H26").Value = 0
$ws.Range("J26").Value = 0
$ws.Range("L26").Value = 0
$ws.Range("N26").ClearContents()
$ws.Range("H28").Value = 320.64285
$ws.Range("J28").Value = 1433
$ws.Range("L28").Value = 1433
$ws.Range("N28").Value = -2403
$ws.Range("H33").Value = 575
$ws.Range("I33").Value = 353.11765
$ws.Range("K33").Value = 353.11765
$ws.Range("M33").Value = -124.11765
$ws.Range("H40").Value = 2465.6667
$ws.Range("I40").Value = 2499
$ws.Range("J40").Value = 2449
$ws.Range("K40").Value = 2499
$ws.Range("L40").Value = 2449
$ws.Range("M40").Value = -2324
$ws.Range("N40").Value = -2799
$ws.Range("H70").Value = 151659.66
$ws.Range("J70").Value = 3899.875
$ws.Range("L70").Value = 11699.625
$ws.Range("N70").Value = -12239.625
$ws.Range("H73").Value = 151659.66
$ws.Range("J73").Value = 3899.875
$ws.Range("L73").Value = 11699.625
$ws.Range("N73").Value = -13571.625

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 853.0536
$ws.Range("I2").Value = 803.1163
$ws.Range("J2").Value = 1018.2308
$ws.Range("K2").Value = 803.1163
$ws.Range("L2").Value = 1018.2308
$ws.Range("M2").Value = -690.1163
$ws.Range("N2").Value = -1244.2308
$ws.Range("H33").Value = 270000
$ws.Range("I33").Value = 500000
$ws.Range("K33").Value = 500000
$ws.Range("M33").Value = -499671
$ws.Range("H74").Value = 1246.0834
$ws.Range("I74").Value = 1223
$ws.Range("K74").Value = 1223
$ws.Range("M74").Value = -349
$ws.Range("H77").Value = 1246.0834
$ws.Range("I77").Value = 1223
$ws.Range("K77").Value = 6115
$ws.Range("M77").Value = -1747
$ws.Range("H110").Value = 2473.1765
$ws.Range("I110").Value = 2101.2727
$ws.Range("K110").Value = 2101.2727
$ws.Range("M110").Value = -56.27269999999999
$ws.Range("H116").Value = 853.0536
$ws.Range("I116").Value = 803.1163
$ws.Range("J116").Value = 1018.2308
$ws.Range("K116").Value = 803.1163
$ws.Range("L116").Value = 1018.2308
$ws.Range("M116").Value = 1490.8837
$ws.Range("N116").Value = -5606.2308

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 853.0536
$ws.Range("I3").Value = 803.1163
$ws.Range("J3").Value = 1018.2308
$ws.Range("K3").Value = 803.1163
$ws.Range("L3").Value = 1018.2308
$ws.Range("M3").Value = -689.1163
$ws.Range("N3").Value = -1246.2308

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 200149.8
$ws.Range("I7").Value = 200149.8
$ws.Range("K7").Value = 200149.8
$ws.Range("M7").Value = -200036.8
$ws.Range("H31").Value = 5728.1274
$ws.Range("J31").Value = 6974.5
$ws.Range("L31").Value = 6974.5
$ws.Range("N31").Value = -7564.5
$ws.Range("H34").Value = 5728.1274
$ws.Range("J34").Value = 6974.5
$ws.Range("L34").Value = 6974.5
$ws.Range("N34").Value = -7378.5
$ws.Range("H62").Value = 7369.6875
$ws.Range("I62").Value = 9420.700000000001
$ws.Range("K62").Value = 9420.700000000001
$ws.Range("M62").Value = -8796.700000000001
$ws.Range("H65").Value = 7369.6875
$ws.Range("I65").Value = 9420.700000000001
$ws.Range("K65").Value = 47103.5
$ws.Range("M65").Value = -43983.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H29").Value = 347.2857
$ws.Range("J29").Value = 167
$ws.Range("L29").Value = 501
$ws.Range("N29").Value = -1055
$ws.Range("H122").Value = 1291.931
$ws.Range("I122").Value = 1014.36365
$ws.Range("J122").Value = 1461.5555
$ws.Range("K122").Value = 9129.272849999999
$ws.Range("L122").Value = 13153.9995
$ws.Range("M122").Value = -6679.272849999999
$ws.Range("N122").Value = -18053.9995
$ws.Range("H132").Value = 1106.3529
$ws.Range("I132").Value = 1127.8
$ws.Range("J132").Value = 945.5
$ws.Range("K132").Value = 10150.2
$ws.Range("L132").Value = 8509.5
$ws.Range("M132").Value = -7620.199999999999
$ws.Range("N132").Value = -13569.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H22").Value = 10000
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 10000
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 10000
$ws.Range("M22").ClearContents()
$ws.Range("N22").Value = -11058
$ws.Range("H80").Value = 20899.111
$ws.Range("I80").Value = 15900
$ws.Range("J80").Value = 24080.363
$ws.Range("K80").Value = 15900
$ws.Range("L80").Value = 24080.363
$ws.Range("M80").Value = -14902
$ws.Range("N80").Value = -26076.363
$ws.Range("H83").Value = 20899.111
$ws.Range("I83").Value = 15900
$ws.Range("J83").Value = 24080.363
$ws.Range("K83").Value = 79500
$ws.Range("L83").Value = 120401.815
$ws.Range("M83").Value = -74508
$ws.Range("N83").Value = -130385.815
$ws.Range("H113").Value = 2168.5833
$ws.Range("I113").Value = 2093
$ws.Range("K113").Value = 2093
$ws.Range("M113").Value = 77

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4018.6206
$ws.Range("I40").Value = 3122.2856
$ws.Range("K40").Value = 3122.2856
$ws.Range("M40").Value = -2986.2856
$ws.Range("H93").Value = 4260.727
$ws.Range("I93").Value = 4211.1875
$ws.Range("K93").Value = 4211.1875
$ws.Range("M93").Value = -2963.1875
$ws.Range("H100").Value = 2521.75
$ws.Range("J100").Value = 2521.75
$ws.Range("L100").Value = 2521.75
$ws.Range("N100").Value = -3603.75
$ws.Range("H132").Value = 2690.238
$ws.Range("J132").Value = 3000
$ws.Range("L132").Value = 9000
$ws.Range("N132").Value = -14060

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H8").Value = 0
$ws.Range("J8").Value = 0
$ws.Range("L8").Value = 0
$ws.Range("N8").ClearContents()
$ws.Range("H21").Value = 0
$ws.Range("I21").Value = 0
$ws.Range("K21").Value = 0
$ws.Range("M21").ClearContents()
$ws.Range("H35").Value = 0
$ws.Range("I35").Value = 0
$ws.Range("K35").Value = 0
$ws.Range("M35").ClearContents()
$ws.Range("H122").Value = 13751.272
$ws.Range("I122").Value = 13602.723
$ws.Range("K122").Value = 40808.169
$ws.Range("M122").Value = -38358.169

